# Edit BurnDownChart & LogBook
#
# - Sheet1 ("LogBook") row edits for the "Benchmark" story tasks:
#     * BootLoop (row 14)        E:J  4 -> 6
#     * Masalah Umum (row 15)    E:J  7 -> 4
#     * Apa itu Benchmark (r16)  E:J  3 -> 1
#     * task in row 17 renamed "Benchmark Apps" -> "Kegunaan Benchmark", E:J 6 -> 1
#     * task in row 18 renamed "Re-design tampilan" -> "Benchmark Apps" (value unchanged)
#   The E21:Q21 (ideal burndown) and E22:Q22 (actual burndown) rows are formula driven
#   off of this table, so they (and the BurnDownChart that reads them) recompute
#   automatically.
# - Selection/view: active cell moves from H21 to H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- LogBook task-hour updates (columns E:J share one value per row here) ---
$ws.Range("E14:J14").Value = 6
$ws.Range("E15:J15").Value = 4
$ws.Range("E16:J16").Value = 1
$ws.Range("E17:J17").Value = 1

# --- Task name swap in column D ---
$ws.Range("D17").Value = "Kegunaan Benchmark"
$ws.Range("D18").Value = "Benchmark Apps"

# --- Update the active cell / selection shown when the workbook is reopened ---
$ws.Activate()
$ws.Range("H19").Select() | Out-Null
